# Close #50 - Loan creation log implementation
# - Add PrestamosLog (loan log) and RetirosLog (withdrawal log) sheets
# - Update Tools with two new tool rows
# - Drop the "Expiracion" column from Loans
# - Refresh sheet view state (zoom/selection/active tab)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Tools sheet - append two new tool rows
# ---------------------------------------------------------------
$toolsWs = $wb.Worksheets.Item("Tools")
$toolsWs.Cells.Item(5, 1).Value = 2
$toolsWs.Cells.Item(5, 2).Value = "Taladro #2"
$toolsWs.Cells.Item(5, 3).Value = 1
$toolsWs.Cells.Item(6, 1).Value = 3
$toolsWs.Cells.Item(6, 2).Value = "Destornillador #2"
$toolsWs.Cells.Item(6, 3).Value = 2

# ---------------------------------------------------------------
# 2. Loans sheet - drop the "Expiracion" column (column E)
# ---------------------------------------------------------------
$loansWs = $wb.Worksheets.Item("Loans")
$loansWs.Range("E1:E4").ClearContents()

# ---------------------------------------------------------------
# 3. Add the PrestamosLog sheet right after Loans
# ---------------------------------------------------------------
$prestamosWs = $wb.Worksheets.Add($null, $loansWs)
$prestamosWs.Name = "PrestamosLog"

$prestamosWs.Cells.Item(1, 1).Value = "Numero"
$prestamosWs.Cells.Item(1, 2).Value = "Elemento"
$prestamosWs.Cells.Item(1, 3).Value = "Trabajador"
$prestamosWs.Cells.Item(1, 4).Value = "Accion"
$prestamosWs.Cells.Item(1, 5).Value = "Fecha"

$prestamosWs.Cells.Item(2, 1).Value = 0
$prestamosWs.Cells.Item(2, 2).Value = "Taladro #1"
$prestamosWs.Cells.Item(2, 3).Value = "Goku"
$prestamosWs.Cells.Item(2, 4).Value = "Apertura"
$prestamosWs.Cells.Item(2, 5).Value = "Tue May 22 16:42:08 GMT-03:00 2018"

$prestamosWs.Cells.Item(3, 1).Value = 1
$prestamosWs.Cells.Item(3, 2).Value = "Destornillador #1"
$prestamosWs.Cells.Item(3, 3).Value = "Vegeta"
$prestamosWs.Cells.Item(3, 4).Value = "Apertura"
$prestamosWs.Cells.Item(3, 5).Value = "Tue May 22 16:42:08 GMT-03:00 2018"

$prestamosWs.Cells.Item(4, 1).Value = 2
$prestamosWs.Cells.Item(4, 2).Value = "Martillo #1"
$prestamosWs.Cells.Item(4, 3).Value = "Vegeta"
$prestamosWs.Cells.Item(4, 4).Value = "Apertura"
$prestamosWs.Cells.Item(4, 5).Value = "Tue May 22 16:42:08 GMT-03:00 2018"

# These 4 log rows record open/close actions on "Taladro #2" - stored as text, not numbers
$prestamosWs.Range("A5:A8").NumberFormat = "@"

$prestamosWs.Cells.Item(5, 1).Value = "3"
$prestamosWs.Cells.Item(5, 2).Value = "Taladro #2"
$prestamosWs.Cells.Item(5, 3).Value = "Goku"
$prestamosWs.Cells.Item(5, 4).Value = "Apertura"
$prestamosWs.Cells.Item(5, 5).Value = "Thu May 31 22:26:34 ART 2018"

$prestamosWs.Cells.Item(6, 1).Value = "3"
$prestamosWs.Cells.Item(6, 2).Value = "Taladro #2"
$prestamosWs.Cells.Item(6, 3).Value = "Goku"
$prestamosWs.Cells.Item(6, 4).Value = "Apertura"
$prestamosWs.Cells.Item(6, 5).Value = "Thu May 31 22:26:35 ART 2018"

$prestamosWs.Cells.Item(7, 1).Value = "3"
$prestamosWs.Cells.Item(7, 2).Value = "Taladro #2"
$prestamosWs.Cells.Item(7, 3).Value = "Goku"
$prestamosWs.Cells.Item(7, 4).Value = "Cierre"
$prestamosWs.Cells.Item(7, 5).Value = "Thu May 31 22:26:35 ART 2018"

$prestamosWs.Cells.Item(8, 1).Value = "3"
$prestamosWs.Cells.Item(8, 2).Value = "Taladro #2"
$prestamosWs.Cells.Item(8, 3).Value = "Goku"
$prestamosWs.Cells.Item(8, 4).Value = "Apertura"
$prestamosWs.Cells.Item(8, 5).Value = "Thu May 31 22:26:35 ART 2018"

$prestamosWs.Columns.Item(2).ColumnWidth = 14.2550
$prestamosWs.Columns.Item(3).ColumnWidth = 10.6300
$prestamosWs.Columns.Item(5).ColumnWidth = 31.2550
$prestamosWs.PageSetup.PaperSize = 9
$prestamosWs.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# 4. Add the RetirosLog sheet right after PrestamosLog
# ---------------------------------------------------------------
$retirosWs = $wb.Worksheets.Add($null, $prestamosWs)
$retirosWs.Name = "RetirosLog"

$retirosWs.Cells.Item(1, 1).Value = "Elemento"
$retirosWs.Cells.Item(1, 2).Value = "Trabajador"
$retirosWs.Cells.Item(1, 3).Value = "Cantidad"
$retirosWs.Cells.Item(1, 4).Value = "Fecha"

$retirosWs.Columns.Item(1).ColumnWidth = 8.0050
$retirosWs.PageSetup.PaperSize = 9
$retirosWs.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# 5. Refresh view state (zoom + selection) across sheets
# ---------------------------------------------------------------
$toolsWs.Range("C6").Select()
$excel.ActiveWindow.Zoom = 400

$suppliesWs = $wb.Worksheets.Item("Supplies")
$suppliesWs.Range("C1:D1").Select()
$excel.ActiveWindow.Zoom = 399

$borrowersWs = $wb.Worksheets.Item("Borrowers")
$borrowersWs.Range("A4:XFD5").Select()
$excel.ActiveWindow.Zoom = 400

$loansWs.Range("A2").Select()
$excel.ActiveWindow.Zoom = 356

$retirosWs.Range("C2").Select()
$excel.ActiveWindow.Zoom = 400

$tipoWs = $wb.Worksheets.Item("TipoHerramienta")
$tipoWs.Range("B7").Select()
$excel.ActiveWindow.Zoom = 400

# PrestamosLog is the sheet that ends up active/selected
$prestamosWs.Range("D6").Select()
$excel.ActiveWindow.Zoom = 385
